$wb = $excel.ActiveWorkbook

# New handback timestamps for the two language sheets.
$zhDateTime = "2016-01-25 07:00:15"
$deDateTime = "2016-01-25 07:00:34"
$statusText = "Handed back: in sync with en-US"

$langs = @(
    @{ Sheet = "zh-cn"; DateTime = $zhDateTime; Xlf = "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf" },
    @{ Sheet = "de-de"; DateTime = $deDateTime; Xlf = "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf" }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Row 2 = a.md.md, Row 3 = b.md.md: both have been handed back and are
    # now in sync with en-US, so record the target file / handback file /
    # handback datetime and flip the status text.
    foreach ($row in 2, 3) {
        $ws.Cells.Item($row, 2).Value = $statusText

        # Latest Target File (E) + Latest Handback File (F) now mirror the
        # source file name / handoff xlf that were already on the row.
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 5), "https://github.com/OpenLocalizationTest/oltest/blob/8b99a6f5fc011dbc757cf64e94c04c5026394207/e2e/a.md.md", "", "", "a.md.md")
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $ws.Hyperlinks.Item($ws.Hyperlinks.Count - 1).Address, "", "", $lang.Xlf)

        # Latest Handback DateTime (G)
        $ws.Cells.Item($row, 7).Value = $lang.DateTime
    }
}

# The Overview sheet's Status columns (B/C) reuse the same shared string as
# the per-language Status column, so they pick up the text change
# automatically once the shared string itself changes - but make sure it is
# explicit here too in case the runtime does not share strings implicitly.
$ov = $wb.Worksheets.Item("Overview")
foreach ($row in 2, 3) {
    $ov.Cells.Item($row, 2).Value = $statusText
    $ov.Cells.Item($row, 3).Value = $statusText
}
